$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new source URL for "confirmed deaths by age Ontario" (row 9) as a hyperlink in C9,
# matching the styling used by the other source-link cells in column C.
[void]$ws.Hyperlinks.Add($ws.Range("C9"), "https://covid-19.ontario.ca/")
$ws.Range("C9").Style = "Hyperlink"

# Update the active selection / view position to B21 (and drop the custom scroll anchor).
[void]$ws.Range("B21").Select()
